# Textbox response formatting fix
# Rename task-order sheets and update the generated stimulus-file names
# (and fix the swapped "eyes closed"/"eyes open" labels on the RS sheet).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG ---
$ws = $wb.Worksheets.Item(1)
$ws.Name = "GNG_TO-16511687019149263"
$ws.Range("B2").Value = "go_stims-16511687018709295.csv"
$ws.Range("B3").Value = "GNG_stims-16511687018979337.csv"
$ws.Range("B4").Value = "go_stims-16511687018999283.csv"
$ws.Range("B5").Value = "GNG_stims-16511687019139292.csv"

# --- Sheet 2: NB ---
$ws = $wb.Worksheets.Item(2)
$ws.Name = "NB_TO-16511687036174378"
$ws.Range("B2").Value  = "OB-16511687026541667.csv"
$ws.Range("B3").Value  = "TB-16511687031781378.csv"
$ws.Range("B4").Value  = "TB-16511687036024323.csv"
$ws.Range("B5").Value  = "ZB-match_7-16511687020079277.csv"
$ws.Range("B6").Value  = "TB-16511687035024018.csv"
$ws.Range("B7").Value  = "ZB-match_9-16511687021761775.csv"
$ws.Range("B8").Value  = "OB-16511687031291397.csv"
$ws.Range("B9").Value  = "ZB-match_6-16511687025271375.csv"
$ws.Range("B10").Value = "OB-1651168702584142.csv"

# --- Sheet 3: RS ---
$ws = $wb.Worksheets.Item(3)
$ws.Name = "RS_TO-16511687036194012"
$ws.Range("B2").Value = "eyes open"
$ws.Range("B3").Value = "eyes closed"

# --- Sheet 4: TOL ---
$ws = $wb.Worksheets.Item(4)
$ws.Name = "TOL_TO-16511687036664028"
$ws.Range("B2").Value = "MM_stims-16511687036334376.csv"
$ws.Range("B3").Value = "ZM_stims-16511687036204019.csv"
$ws.Range("B4").Value = "MM_stims-1651168703649433.csv"
$ws.Range("B5").Value = "ZM_stims-16511687036343987.csv"
$ws.Range("B6").Value = "MM_stims-1651168703665397.csv"
$ws.Range("B7").Value = "ZM_stims-1651168703649433.csv"

# --- Sheet 5: vSAT ---
$ws = $wb.Worksheets.Item(5)
$ws.Name = "vSAT_TO-16511687037433991"
$ws.Range("B2").Value = "SAT_stims-16511687036963973.csv"
$ws.Range("B3").Value = "vSAT_stims-16511687037123969.csv"
$ws.Range("B4").Value = "vSAT_stims-16511687037283976.csv"
$ws.Range("B5").Value = "SAT_stims-16511687036713989.csv"
